# Regen s_val data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 0, 11.16088251738829),
    @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 1.271902929317955),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 1, 13.71653804550039),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 1, 11.945164432584),
    @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1, 1.649481363816475),
    @(0.0000000008413942875762359, 0.002777888934908601, 0.8054896365839992, 616238.5361209477, 0, 616239.344388474),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182),
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797),
    @(1.459612070389937, 1.667794583268128, 26.21740644021617, 8.660232485948974, 0, 38.00504557982321),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759),
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 4.429675500412797)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
